$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Occurrence" column (E) for rows 2-12 shares one string value; update it
# from "9, 13, 36" to "9, 13, 36, 40" (row 13 keeps its own "9, 13" value).
for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    if ($cell.Value2 -eq "9, 13, 36") {
        $cell.Value = "9, 13, 36, 40"
    }
}

# Match the saved view state: selection moved from E13 to E12, and the
# window was scrolled up one row (top-left visible cell A10 -> A9).
$ws.Activate()
$ws.Range("E12").Select()

$win = $excel.ActiveWindow
$win.ScrollRow = 9
$win.ScrollColumn = 1
